# Daily attendance processing - 2025-12-30 18:09:29
# Normalize the "Recorded By" (column G) entries so that the "System"
# recorder is always listed first among the recorders for a session,
# preserving the relative order of the remaining recorders. When no
# "System" entry is present but two recorders are listed, their order
# is swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cellVal = $cell.Text

    if ([string]::IsNullOrEmpty($cellVal)) {
        continue
    }

    $rawParts = $cellVal.Split(",")
    if ($rawParts.Count -lt 2) {
        continue
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $sysIdx = [array]::IndexOf($parts, "System")

    $newParts = $null
    if ($sysIdx -ge 0) {
        $rest = @()
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($i -ne $sysIdx) {
                $rest += $parts[$i]
            }
        }
        $newParts = @("System") + $rest
    } elseif ($parts.Count -eq 2) {
        $newParts = @($parts[1], $parts[0])
    }

    if ($newParts -ne $null) {
        $newVal = $newParts -join ", "
        if ($newVal -ne $cellVal) {
            $cell.Value = $newVal
        }
    }
}
